$d = $word.ActiveDocument
$nbsp = [char]160

# 1. "Allowed-with list version" table: 1.1 -> 1.0e
$r1 = $d.Content.Find.Execute("1.1", $true, $false, $false, $false, $false, $true, 1, $false, "1.0e", 1)

# 2. "Date" table: 16 August 2023<nbsp> -> 15 February 2024
$dateOld = "16 August 2023" + $nbsp
$r2 = $d.Content.Find.Execute($dateOld, $true, $false, $false, $false, $false, $true, 1, $false, "15 February 2024", 1)


# 3. "Latest definitive cPP version" row: 1.1 (16 August 2023) -> 1.0e (15 February 2024)
$r3 = $d.Content.Find.Execute("1.1 (16 August 2023)", $true, $false, $false, $false, $false, $true, 1, $false, "1.0e (15 February 2024)", 1)


# 4. "Latest definitive SD version" row: 1.1 (16 August 2023) -> 1.0e (15 February 2024)
$r4 = $d.Content.Find.Execute("1.1 (16 August 2023)", $true, $false, $false, $false, $false, $true, 1, $false, "1.0e (15 February 2024)", 1)


# 5. PP-Module for Server Applications, Object version row: reduce leading double
#    space to a single space and update the version/date text.
$r5 = $d.Content.Find.Execute("  1.1 (16 August 2023)", $true, $false, $false, $false, $false, $true, 1, $false, " 1.0e (15 February 2024)", 1)


# 6. PP-Module for Agent Applications, Object version row: same text update,
#    but the leading double space is left untouched.
$r6 = $d.Content.Find.Execute("1.1 (16 August 2023)", $true, $false, $false, $false, $false, $true, 1, $false, "1.0e (15 February 2024)", 1)

Write-Host "step1:" $r1 " step2:" $r2 " step3:" $r3 " step4:" $r4 " step5:" $r5 " step6:" $r6
